$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 377; this shifts rows 377:457 down to 378:458
# and extends the used range to A1:R458 (matches the diff's dimension change).
$ws.Rows(377).Insert()

# Populate the newly inserted row 377 with the new daily price record.
$ws.Range("A377").Value = 4
$ws.Range("B377").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C377").Value = "Los Lagos"
$ws.Range("D377").Value = 45173
$ws.Range("E377").Value = 10
$ws.Range("F377").Value = 100112032
$ws.Range("G377").Value = "Zapallo italiano"
$ws.Range("H377").Value = "Sin especificar"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 70
$ws.Range("K377").Value = 20000
$ws.Range("L377").Value = 20000
$ws.Range("M377").Value = 20000
$ws.Range("N377").Value = "$/caja 50 unidades"
$ws.Range("O377").Value = "Región de Arica y Parinacota"
$ws.Range("P377").Value = 400
$ws.Range("Q377").Value = 50
$ws.Range("R377").Value = "Hortaliza"
